# login fonctionel avec traitement des signIn/Up
# Replays the sign-up pseudo testing that happened on the "Personnes" sheet:
# several candidate pseudos were typed into B2 and B3 before the final
# values ("bernard" / "bernarda") were kept.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

# Candidate pseudos tried out in B2 while testing the sign up feature
$ws.Range("B2").Value = "test_ce_pseudo"
$ws.Range("B2").Value = "bernardoo"
$ws.Range("B2").Value = "natafa"
$ws.Range("B2").Value = "bernard"

# Candidate pseudos tried out in B3
$ws.Range("B3").Value = "bernardo"
$ws.Range("B3").Value = "bernarda"
